# "Thêm giao diện Ban Đào tạo" — add a new account row for the "Ban Đào Tạo"
# (Training Department) role, mirroring the existing "BanNhanSu" (super user)
# row already on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 8: UserName / Password / UserType, following the same pattern as
# row 7 (A7="BanNhanSu", B7=12345, C7="super").
$ws.Range("A8").Value = "BanDaoTao"
$ws.Range("B8").Value = 12345
$ws.Range("C8").Value = "super"

# Column A was given an explicit width in the authored workbook.
$ws.Columns("A").ColumnWidth = 10.5

# Move the active selection to D8, matching where the author's cursor ended
# up after entering the new row.
$ws.Range("D8").Select() | Out-Null
